$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the born/relive position ("RelivePos") for the first scene row
# from "0,0,0" to "20,0,-137".
$ws.Range("E2").Value = "20,0,-137"
